# syncing prior to aws migration
# Add a new "Fudge Factor" column (O1 header / O2 value) to the Variables sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# New header cell, matching the style of the existing header row (e.g. M1).
$ws.Range("O1").Value = "Fudge Factor"
$ws.Range("M1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data cell under the header.
$ws.Range("O2").Value = 3

# Update the active selection to reflect where editing left off.
[void]$ws.Range("N4").Select()

# Reflect the updated window layout recorded for this edit.
$w = $excel.ActiveWindow
$w.Left = 390
$w.Top = 390
$w.Width = 28800
$w.Height = 14430
